# This edit reorders the content of several rows in the "Artfynd" sheet.
# Each physical row keeps its position/formatting, but the field values that
# used to live in one row now belong to a different row within the same
# small cluster of rows (a straightforward swap for pairs of rows, and a
# three-way rotation for a few clusters of three rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: Excel's COM layer auto-detects strings that look like full dates
# (e.g. "2026-01-18") and silently converts them to date values when such a
# value is assigned to a cell/range. The source cells are plain text, so we
# guard against that by forcing a leading apostrophe (text qualifier) on any
# value that matches a yyyy-mm-dd pattern before writing it back.
function Protect-DateStrings($arr) {
    $cols = $arr.GetUpperBound(1)
    for ($j = 1; $j -le $cols; $j++) {
        $val = $arr[1, $j]
        if ($val -is [string] -and $val -match '^\d{4}-\d{2}-\d{2}$') {
            $arr[1, $j] = "'" + $val
        }
    }
    return $arr
}

# Grab the full A:AY content of every row involved in the reshuffle before
# writing anything back, so the three-way rotations do not clobber data that
# is still needed for a later assignment.
$row2 = $ws.Range("A2:AY2").Value()
$row3 = $ws.Range("A3:AY3").Value()
$row4 = $ws.Range("A4:AY4").Value()
$row5 = $ws.Range("A5:AY5").Value()
$row7 = $ws.Range("A7:AY7").Value()
$row8 = $ws.Range("A8:AY8").Value()
$row9 = $ws.Range("A9:AY9").Value()
$row19 = $ws.Range("A19:AY19").Value()
$row20 = $ws.Range("A20:AY20").Value()
$row21 = $ws.Range("A21:AY21").Value()
$row22 = $ws.Range("A22:AY22").Value()
$row23 = $ws.Range("A23:AY23").Value()
$row25 = $ws.Range("A25:AY25").Value()
$row26 = $ws.Range("A26:AY26").Value()
$row27 = $ws.Range("A27:AY27").Value()
$row28 = $ws.Range("A28:AY28").Value()
$row34 = $ws.Range("A34:AY34").Value()
$row35 = $ws.Range("A35:AY35").Value()
$row36 = $ws.Range("A36:AY36").Value()
$row42 = $ws.Range("A42:AY42").Value()
$row43 = $ws.Range("A43:AY43").Value()
$row44 = $ws.Range("A44:AY44").Value()

# Write each row's new content into place.
# Rows 2, 3 swap
$ws.Range("A2:AY2").Value = (Protect-DateStrings $row3)
$ws.Range("A3:AY3").Value = (Protect-DateStrings $row2)

# Rows 4, 5 swap
$ws.Range("A4:AY4").Value = (Protect-DateStrings $row5)
$ws.Range("A5:AY5").Value = (Protect-DateStrings $row4)

# Rows 7, 8, 9 rotate (new7=old9, new8=old7, new9=old8)
$ws.Range("A7:AY7").Value = (Protect-DateStrings $row9)
$ws.Range("A9:AY9").Value = (Protect-DateStrings $row8)
$ws.Range("A8:AY8").Value = (Protect-DateStrings $row7)

# Rows 19, 20, 21 rotate (new19=old20, new20=old21, new21=old19)
$ws.Range("A19:AY19").Value = (Protect-DateStrings $row20)
$ws.Range("A20:AY20").Value = (Protect-DateStrings $row21)
$ws.Range("A21:AY21").Value = (Protect-DateStrings $row19)

# Rows 22, 23 swap
$ws.Range("A22:AY22").Value = (Protect-DateStrings $row23)
$ws.Range("A23:AY23").Value = (Protect-DateStrings $row22)

# Rows 25, 26 swap
$ws.Range("A25:AY25").Value = (Protect-DateStrings $row26)
$ws.Range("A26:AY26").Value = (Protect-DateStrings $row25)

# Rows 27, 28 swap
$ws.Range("A27:AY27").Value = (Protect-DateStrings $row28)
$ws.Range("A28:AY28").Value = (Protect-DateStrings $row27)

# Rows 34, 35, 36 rotate (new34=old36, new35=old34, new36=old35)
$ws.Range("A34:AY34").Value = (Protect-DateStrings $row36)
$ws.Range("A35:AY35").Value = (Protect-DateStrings $row34)
$ws.Range("A36:AY36").Value = (Protect-DateStrings $row35)

# Rows 42, 43, 44 rotate (new42=old43, new43=old44, new44=old42)
$ws.Range("A42:AY42").Value = (Protect-DateStrings $row43)
$ws.Range("A43:AY43").Value = (Protect-DateStrings $row44)
$ws.Range("A44:AY44").Value = (Protect-DateStrings $row42)
